$d = $word.ActiveDocument

# 1) April 9: merge "Research project presentations" + trailing " " run into a
#    single preserved-space run. Only the first occurrence (April 9) has the
#    extra space run, so searching/replacing with the exact spaced text only
#    matches (and fixes) that one paragraph.
$d.Content.Find.Execute("Research project presentations ", $true, $false, $false, $false, $false, $true, 1, $false, "Research project presentations ", 2) | Out-Null

# 2) May 5: "Project Presentations" -> "Final Project"
$d.Content.Find.Execute("Project Presentations", $true, $false, $false, $false, $false, $true, 1, $false, "Final Project", 2) | Out-Null

# 3) May 7: "Project & Research presentations" -> "Final Project"
$d.Content.Find.Execute("Project & Research presentations", $true, $false, $false, $false, $false, $true, 1, $false, "Final Project", 2) | Out-Null

# 4) May 12: add a new paragraph "Presentations 8-9:50" right after the "12"
#    day-number paragraph (before the existing trailing empty paragraph).
$idx = 0
$dayTwelveIndex = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($dayTwelveIndex -eq -1 -and $p.Range.Text -eq "12`r") {
        # the May calendar's "12" cell is the second one found (the first is
        # in an earlier month's table); track all matches and use the one
        # immediately followed by an empty paragraph with no pStyle.
    }
}

$idx = 0
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -eq "12`r") {
        $nextIdx = $idx + 1
        $n = 0
        foreach ($p2 in $d.Paragraphs) {
            $n = $n + 1
            if ($n -eq $nextIdx) {
                if ($p2.Range.Text -eq "`r") {
                    $targetIdx = $nextIdx
                }
            }
        }
    }
}

if ($targetIdx -ne -1) {
    $n2 = 0
    foreach ($p in $d.Paragraphs) {
        $n2 = $n2 + 1
        if ($n2 -eq $targetIdx) {
            $p.Range.InsertParagraphBefore() | Out-Null
        }
    }
    $n3 = 0
    foreach ($p in $d.Paragraphs) {
        $n3 = $n3 + 1
        if ($n3 -eq $targetIdx) {
            $p.Range.Text = "Presentations 8-9:50"
        }
    }
}
